$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 84.97741527679979
$ws.Range("B3").Value = 88.51599751506561
$ws.Range("B4").Value = 91.32406781182483
$ws.Range("I5").Value = 95.9088214900177
$ws.Range("I6").Value = 95.92836530352579
$ws.Range("I7").Value = 95.93481981522186
$ws.Range("C8").Value = 98.55740384435178
$ws.Range("C9").Value = 97.24656813527551
$ws.Range("C10").Value = 98.19069397979416
$ws.Range("D11").Value = 99.30371935120795
$ws.Range("D12").Value = 99.19919535414697
$ws.Range("D13").Value = 99.2516874970003
$ws.Range("E14").Value = 98.71552424447512
$ws.Range("E15").Value = 98.7746377656615
$ws.Range("E16").Value = 98.73649736466054
$ws.Range("F17").Value = 98.20670302790766
$ws.Range("F18").Value = 98.27512025665715
$ws.Range("F19").Value = 98.19310947726699
$ws.Range("G20").Value = 97.71525622620572
$ws.Range("G21").Value = 97.80594482642393
$ws.Range("G22").Value = 97.76030576675177
$ws.Range("H23").Value = 97.29158969819912
$ws.Range("H24").Value = 97.28853691751796
$ws.Range("H25").Value = 97.26979382750832
$ws.Range("B26").Value = 90.11029093408622
$ws.Range("B27").Value = 93.16418304623068
$ws.Range("I28").Value = 95.82303974401631
$ws.Range("I29").Value = 95.85106016982516
$ws.Range("C30").Value = 97.81416658068642
$ws.Range("C31").Value = 98.02094649342695
$ws.Range("D32").Value = 99.26498696841793
$ws.Range("D33").Value = 99.25356781483723
$ws.Range("E34").Value = 98.76204864640791
$ws.Range("E35").Value = 98.72347453215579
$ws.Range("F36").Value = 98.23945858993132
$ws.Range("F37").Value = 98.30902888227904
$ws.Range("G38").Value = 97.71099718609373
$ws.Range("G39").Value = 97.71420039556003
$ws.Range("H40").Value = 97.304207780216
$ws.Range("H41").Value = 97.27949885887465
$ws.Range("B42").Value = 91.61354032889869
$ws.Range("B43").Value = 94.27069900561777
$ws.Range("I44").Value = 95.8541935890678
$ws.Range("I45").Value = 95.93244963558024
$ws.Range("C46").Value = 98.56857693231979
$ws.Range("C47").Value = 98.26408196101325
$ws.Range("D48").Value = 99.27163357741426
$ws.Range("D49").Value = 99.21715802886678
$ws.Range("E50").Value = 98.71466997249088
$ws.Range("E51").Value = 98.65424833470159
$ws.Range("F52").Value = 98.2767620937244
$ws.Range("F53").Value = 98.14483920762531
$ws.Range("G54").Value = 97.70611993284349
$ws.Range("G55").Value = 97.75179399754147
$ws.Range("H56").Value = 97.2345654022181
$ws.Range("H57").Value = 97.15811035282969
